$wb = $excel.ActiveWorkbook

# ---- Sheet "Главные" (Main) ----
$ws2 = $wb.Worksheets.Item("Главные")

$ws2.Range("C4").Value = 18
$ws2.Range("D4").Value = 296
$ws2.Range("F4").Value = 171
$ws2.Range("G4").Value = 16.44
$ws2.Range("H4").Value = 6.94
$ws2.Range("I4").Value = 9.5
$ws2.Range("K4").Value = 73
$ws2.Range("C8").Value = 21
$ws2.Range("D8").Value = 394
$ws2.Range("E8").Value = 195
$ws2.Range("F8").Value = 199
$ws2.Range("G8").Value = 18.76
$ws2.Range("H8").Value = 9.289999999999999
$ws2.Range("I8").Value = 9.48
$ws2.Range("J8").Value = 90
$ws2.Range("K8").Value = 92
$ws2.Range("C9").Value = 24
$ws2.Range("D9").Value = 374
$ws2.Range("E9").Value = 201
$ws2.Range("F9").Value = 173
$ws2.Range("G9").Value = 15.58
$ws2.Range("H9").Value = 8.380000000000001
$ws2.Range("I9").Value = 7.21
$ws2.Range("J9").Value = 98
$ws2.Range("K9").Value = 84
$ws2.Range("C10").Value = 16
$ws2.Range("D10").Value = 292
$ws2.Range("E10").Value = 148
$ws2.Range("F10").Value = 144
$ws2.Range("G10").Value = 18.25
$ws2.Range("H10").Value = 9.25
$ws2.Range("I10").Value = 9
$ws2.Range("J10").Value = 74
$ws2.Range("K10").Value = 62
$ws2.Range("C14").Value = 15
$ws2.Range("D14").Value = 181
$ws2.Range("E14").Value = 96
$ws2.Range("F14").Value = 85
$ws2.Range("G14").Value = 12.07
$ws2.Range("H14").Value = 6.4
$ws2.Range("I14").Value = 5.67
$ws2.Range("J14").Value = 48
$ws2.Range("K14").Value = 40
$ws2.Range("C15").Value = 16
$ws2.Range("D15").Value = 214
$ws2.Range("F15").Value = 124
$ws2.Range("G15").Value = 13.38
$ws2.Range("H15").Value = 5.63
$ws2.Range("I15").Value = 7.75
$ws2.Range("K15").Value = 62
$ws2.Range("C16").Value = 23
$ws2.Range("D16").Value = 456
$ws2.Range("E16").Value = 223
$ws2.Range("F16").Value = 233
$ws2.Range("G16").Value = 19.83
$ws2.Range("H16").Value = 9.699999999999999
$ws2.Range("I16").Value = 10.13
$ws2.Range("J16").Value = 84
$ws2.Range("K16").Value = 84
$ws2.Range("C18").Value = 23
$ws2.Range("D18").Value = 338
$ws2.Range("E18").Value = 158
$ws2.Range("F18").Value = 180
$ws2.Range("G18").Value = 14.7
$ws2.Range("H18").Value = 6.87
$ws2.Range("I18").Value = 7.83
$ws2.Range("J18").Value = 69
$ws2.Range("K18").Value = 85
$ws2.Range("C19").Value = 18
$ws2.Range("D19").Value = 328
$ws2.Range("E19").Value = 154
$ws2.Range("F19").Value = 174
$ws2.Range("G19").Value = 18.22
$ws2.Range("H19").Value = 8.56
$ws2.Range("I19").Value = 9.67
$ws2.Range("J19").Value = 72
$ws2.Range("K19").Value = 72
$ws2.Range("C20").Value = 22
$ws2.Range("D20").Value = 377
$ws2.Range("E20").Value = 158
$ws2.Range("F20").Value = 219
$ws2.Range("G20").Value = 17.14
$ws2.Range("H20").Value = 7.18
$ws2.Range("I20").Value = 9.949999999999999
$ws2.Range("J20").Value = 74
$ws2.Range("K20").Value = 82
$ws2.Range("C21").Value = 19
$ws2.Range("D21").Value = 277
$ws2.Range("E21").Value = 122
$ws2.Range("F21").Value = 155
$ws2.Range("G21").Value = 14.58
$ws2.Range("H21").Value = 6.42
$ws2.Range("I21").Value = 8.16
$ws2.Range("J21").Value = 51
$ws2.Range("K21").Value = 65
$ws2.Range("C22").Value = 17
$ws2.Range("D22").Value = 322
$ws2.Range("E22").Value = 128
$ws2.Range("F22").Value = 194
$ws2.Range("G22").Value = 18.94
$ws2.Range("H22").Value = 7.53
$ws2.Range("I22").Value = 11.41
$ws2.Range("J22").Value = 64
$ws2.Range("K22").Value = 67

# Update as_of_utc timestamps for rows 2-26
for ($r = 2; $r -le 26; $r++) {
    $ws2.Range("AA$r").Value = "2025-11-11 07:06:04"
}

# ---- Sheet "Линейные" (Linear) ----
$ws3 = $wb.Worksheets.Item("Линейные")

$ws3.Range("C2").Value = 14
$ws3.Range("D2").Value = 274
$ws3.Range("E2").Value = 117
$ws3.Range("F2").Value = 157
$ws3.Range("G2").Value = 19.57
$ws3.Range("H2").Value = 8.359999999999999
$ws3.Range("I2").Value = 11.21
$ws3.Range("J2").Value = 51
$ws3.Range("K2").Value = 56
$ws3.Range("L2").Value = 1
$ws3.Range("M2").Value = 1
$ws3.Range("N2").Value = 1
$ws3.Range("V2").Value = 6
$ws3.Range("C3").Value = 22
$ws3.Range("D3").Value = 322
$ws3.Range("E3").Value = 156
$ws3.Range("F3").Value = 166
$ws3.Range("G3").Value = 14.64
$ws3.Range("H3").Value = 7.09
$ws3.Range("I3").Value = 7.55
$ws3.Range("J3").Value = 78
$ws3.Range("K3").Value = 68
$ws3.Range("C4").Value = 12
$ws3.Range("D4").Value = 192
$ws3.Range("E4").Value = 84
$ws3.Range("F4").Value = 108
$ws3.Range("G4").Value = 16
$ws3.Range("H4").Value = 7
$ws3.Range("I4").Value = 9
$ws3.Range("J4").Value = 42
$ws3.Range("K4").Value = 44
$ws3.Range("C5").Value = 10
$ws3.Range("D5").Value = 148
$ws3.Range("E5").Value = 78
$ws3.Range("F5").Value = 70
$ws3.Range("G5").Value = 14.8
$ws3.Range("H5").Value = 7.8
$ws3.Range("I5").Value = 7
$ws3.Range("J5").Value = 39
$ws3.Range("K5").Value = 35
$ws3.Range("C6").Value = 14
$ws3.Range("D6").Value = 271
$ws3.Range("E6").Value = 121
$ws3.Range("F6").Value = 150
$ws3.Range("G6").Value = 19.36
$ws3.Range("H6").Value = 8.640000000000001
$ws3.Range("I6").Value = 10.71
$ws3.Range("J6").Value = 53
$ws3.Range("K6").Value = 70
$ws3.Range("L6").Value = 1
$ws3.Range("M6").Value = 2
$ws3.Range("N6").Value = 1
$ws3.Range("V6").Value = 6
$ws3.Range("C7").Value = 14
$ws3.Range("D7").Value = 223
$ws3.Range("E7").Value = 76
$ws3.Range("F7").Value = 147
$ws3.Range("G7").Value = 15.93
$ws3.Range("H7").Value = 5.43
$ws3.Range("I7").Value = 10.5
$ws3.Range("J7").Value = 38
$ws3.Range("K7").Value = 46
$ws3.Range("C11").Value = 15
$ws3.Range("D11").Value = 215
$ws3.Range("F11").Value = 115
$ws3.Range("G11").Value = 14.33
$ws3.Range("H11").Value = 6.67
$ws3.Range("I11").Value = 7.67
$ws3.Range("K11").Value = 55
$ws3.Range("C19").Value = 22
$ws3.Range("D19").Value = 393
$ws3.Range("E19").Value = 186
$ws3.Range("F19").Value = 207
$ws3.Range("G19").Value = 17.86
$ws3.Range("H19").Value = 8.449999999999999
$ws3.Range("I19").Value = 9.41
$ws3.Range("J19").Value = 88
$ws3.Range("K19").Value = 91
$ws3.Range("C21").Value = 26
$ws3.Range("D21").Value = 525
$ws3.Range("E21").Value = 217
$ws3.Range("F21").Value = 308
$ws3.Range("G21").Value = 20.19
$ws3.Range("H21").Value = 8.35
$ws3.Range("I21").Value = 11.85
$ws3.Range("J21").Value = 101
$ws3.Range("K21").Value = 124
$ws3.Range("C22").Value = 17
$ws3.Range("D22").Value = 275
$ws3.Range("E22").Value = 134
$ws3.Range("F22").Value = 141
$ws3.Range("G22").Value = 16.18
$ws3.Range("H22").Value = 7.88
$ws3.Range("I22").Value = 8.289999999999999
$ws3.Range("J22").Value = 67
$ws3.Range("K22").Value = 68
$ws3.Range("C24").Value = 24
$ws3.Range("D24").Value = 424
$ws3.Range("F24").Value = 257
$ws3.Range("G24").Value = 17.67
$ws3.Range("H24").Value = 6.96
$ws3.Range("I24").Value = 10.71
$ws3.Range("K24").Value = 101
$ws3.Range("C26").Value = 20
$ws3.Range("D26").Value = 422
$ws3.Range("E26").Value = 179
$ws3.Range("F26").Value = 243
$ws3.Range("G26").Value = 21.1
$ws3.Range("H26").Value = 8.949999999999999
$ws3.Range("I26").Value = 12.15
$ws3.Range("J26").Value = 67
$ws3.Range("K26").Value = 69

# Update as_of_utc timestamps for rows 2-26
for ($r = 2; $r -le 26; $r++) {
    $ws3.Range("AA$r").Value = "2025-11-11 07:06:04"
}
